# Entry Parser tested and Works
# Adds a new test-case row (row 14) to the activity-parser table on Sheet1,
# mirroring row 3 (Input1=Swim, Activity1=Swim) but also populating
# Input2=Swim, and tidies up the leftover formatting on row 13 (the
# "Difficult Competions" row) that no longer needs the extra style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 previously carried a stray alignment style (applied to a
# 6-cell block, including empty cells). Strip that formatting and drop
# the now-pointless empty cells (B13, E13, F13) entirely.
$ws.Range("A13:F13").ClearFormats()
$ws.Range("B13").Clear()
$ws.Range("E13").Clear()
$ws.Range("F13").Clear()

# New test case in row 14: Input 1 = Swim, Input 2 = Swim, Activity 1 = Swim
$ws.Range("A14").Value = "Swim"
$ws.Range("B14").Value = "Swim"
$ws.Range("D14").Value = "Swim"

# The Table1 (A2:F13) grows by one row to include the new entry.
[void]$ws.ListObjects.Item("Table1").Resize($ws.Range("A2:F14"))

# Selection moves on to the next empty row beneath the table, and the
# sheet view no longer pins a frozen/scrolled top-left cell.
[void]$ws.Range("D16").Select()
